$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new feature row (row 14), mirroring the structure of existing rows.
$ws.Cells.Item(14, 1).Value = "Exploding dice"
$ws.Cells.Item(14, 2).Value = "When you roll maximum on a die, keep rolling until you don't"
$ws.Cells.Item(14, 4).Value = "Mike Combs - Store review"

# Update the active selection to match the post-edit workbook state.
$ws.Range("D22").Select()
